$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (phoneNumber), shifting F:M -> G:N
$ws.Columns.Item(6).Insert()

# Old M3 (wrap-text / fill style) lands on N3 after the shift; the new layout
# no longer uses that style, so reset it back to the workbook default.
$ws.Range("N3").Style = "Normal"

# --- Header row (row 1) ---
$ws.Range("F1").Value = "mobile_no_country_code"
$ws.Range("E1").Value = "id_no"

# --- Data row (row 3) ---
$ws.Range("F3").Value = "'+27"

$ws.Range("G1").Value = "mobile_no"

# --- Column widths ---
$ws.Columns.Item(5).ColumnWidth = 6.33203125
$ws.Columns.Item(6).ColumnWidth = 25.109375

# --- Selection ---
$ws.Range("N7").Select()
